$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.561.35"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.505.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.30%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "195.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.71%  "
$ws.Range("E7").Value = "  -1.49%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -5.84%  "
$ws.Range("E10").Value = "  -2.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.07"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.74%  "
$ws.Range("E12").Value = "  -3.38%  "
$ws.Range("E13").Value = "  -1.70%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.060.31"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "598.00"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.44%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.706.70"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("E17").Value = "  -0.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.63"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.536.66"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.79%  "
$ws.Range("E20").Value = "  +2.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.982"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.83%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.90"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "101.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.48%  "
$ws.Range("E25").Value = "  -2.20%  "
$ws.Range("E26").Value = "  +0.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.47"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.49%  "
$ws.Range("E29").Value = "  -3.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.22%  "
$ws.Range("E32").Value = "  -2.04%  "
$ws.Range("E33").Value = "  -2.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.10"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.80%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.21"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.743.23"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.87%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0₃0807"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.53%  "
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.64"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.66%  "
$ws.Range("E40").Value = "  -1.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.03"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "492.26"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.79%  "
$ws.Range("E43").Value = "  -3.67%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0449"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.78%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.81"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.86%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.138"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.23"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.94%  "
$ws.Range("E48").Value = "  +0.14%  "
$ws.Range("E49").Value = "  -4.43%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000243"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.80%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.35"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.11%  "
